$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 changes
$ws.Range("AC3").Value = 6
$ws.Range("AS3").Value = 301

# Row 4 changes
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5

# Row 5 changes
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 3.5
$ws.Range("Q5").Value = 1.98
$ws.Range("R5").Value = 1.88
